# Add season-record columns (Wins / Losses / Ties) to the LAD_2012 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header row (bold, centered, thin border)
# by copying the formatting from the neighboring header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row (2-52) gets the same 2012 Dodgers season record:
# 86 wins, 76 losses, 0 ties.
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}
